# Add the three new character styles (GaNStyle, GaNParagraph, GaNLinks)
# and apply them to the runs that were newly styled in the commit.

$d = $word.ActiveDocument

# wdStyleTypeCharacter = 2
$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608   # wdColor value for RGB(0,0,128) -> OOXML 000080
$gaNLinks.Font.Underline = 1     # wdUnderlineSingle

# --- Apply GaNStyle to every "2022: Daty kampanii ..." run (4 occurrences) ---
$campaignText = "2022: Daty kampanii używające Gwiazdozbiór Pegaza: 8-17 października, 7-16 listopada,"
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End
while ($rng.Find.Execute($campaignText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $rng.End = $d.Content.End
}

# --- Apply GaNParagraph to the "Uczestniczysz w ogólnoświatowym ..." run ---
$paragraphText = "Uczestniczysz w ogólnoświatowym przedsięwzięciu, którego celem jest obserwacja i odnotowanie najsłabszych widocznych gwiazd w celu zmierzenia zanieczyszczenia światłem w danym miejscu. Poprzez zlokalizowanie i obserwację  Gwiazdozbiór Pegaza na nocnym niebie oraz porównanie go do map nieba ludzie z całego świata będą mogli dowiedzieć się jaki wkład światło emitowane przez ich społeczność wnosi do  zanieczyszczenia światłem. To co dodasz do internetowej bazy danych pomoże udokumentować widoczne nocne niebo."
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "Jenika Hollana, CzechGlobe ..." run ---
$linksText = " Jenika Hollana, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3 = $d.Content
$found3 = $rng3.Find.Execute($linksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Style = "GaNLinks"
}
